$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("logs")

# New rows to append: (Id, DescriptionStr, MessageStr, DateSerial)
$newRows = @(
    @(977, "Info", "Tüm öğrenciler listelendi", 45630),
    @(978, "Info", "Tüm öğrenciler listelendi", 45630),
    @(979, "Info", "Tüm öğrenciler listelendi", 45630),
    @(980, "Info", "Tüm öğrenciler listelendi", 45630),
    @(981, "Info", "Tüm öğrenciler listelendi", 45630),
    @(982, "Info", "Tüm öğrenciler listelendi", 45630),
    @(983, "Info", "Tüm öğrenciler listelendi", 45630),
    @(984, "Info", "Tüm öğrenciler listelendi", 45630),
    @(985, "Info", "Tüm öğrenciler listelendi", 45630),
    @(986, "Info", "Tüm öğrenciler listelendi", 45630),
    @(987, "Info", "Tüm öğrenciler listelendi", 45630),
    @(988, "Info", "Tüm öğrenciler listelendi", 45630),
    @(989, "Info", "Tüm öğrenciler listelendi", 45630),
    @(990, "Info", "Tüm departmanlar listelendi", 45630),
    @(991, "Info", "Tüm öğrenciler listelendi", 45631)
)

$startRow = 672
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value = $row[3]
    $cell.NumberFormat = "dd-MM-yyyy"
}
